# Append a new job listing at the top of the "ランサーズ" sheet's data
# table (row 2), pushing the existing rows down by one, and refresh the
# "取得日時" (fetched-at) timestamp on every data row to the new run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-10-20 06:35:12"

# --- 1) Drop every existing hyperlink up front -----------------------
# Hyperlink refs are anchored to absolute cells and do NOT slide down
# automatically when a row is inserted, so if we left them in place the
# row-insert below would leave a stale link (pointing at the old row's
# URL) sitting on the freshly inserted blank row. Clearing them all now
# and re-adding them once every cell holds its final value keeps things
# simple and avoids any duplicate/mismatched links.
$ws.Hyperlinks.Delete()

# --- 2) Insert the new row, shifting rows 2-12 down to 3-13 ----------
$ws.Rows(2).Insert()

# --- 3) Populate the newly inserted row 2 with the new listing -------
$ws.Range("A2").Value = $newTimestamp
$ws.Range("B2").Value = "【効率化】Air Tableでデータ収集と工数管理とスムーズにしたい!"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5416639"
$ws.Range("G2").Value = 388
$ws.Range("H2").Value = "🔥AI,Ai ◆効率化 ◇管理"

# --- 4) Refresh the "取得日時" timestamp on every surviving data row -
for ($r = 3; $r -le 13; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}

# --- 5) Re-apply the Hyperlink style + live hyperlinks on column F ---
for ($r = 2; $r -le 13; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.Style = "Hyperlink"
    $ws.Hyperlinks.Add($cell, $cell.Value)
}
